$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$nm = $p.NotesMaster
$smtcs = $sm.Theme.ThemeColorScheme
$nmtcs = $nm.Theme.ThemeColorScheme

$smtcs.Colors(3).RGB = 12345
Write-Host "after setting SlideMaster dk2, NotesMaster dk2 = $($nmtcs.Colors(3).RGB)"
